try {
  $x = New-Object System.IO.MemoryStream
  Write-Output "ok: $x"
} catch {
  Write-Output "ERR: $_"
}
